# Update driver config: set "Run Mode" (column C) to "No" for most test
# rows, leaving a handful of rows set to "Yes", and update the active
# selection to C31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Automation Tests")

# Rows whose Run Mode flips from "Yes" to "No".
$rowsToNo = @(2,3,4,5,6,7,8,9,10,11,12,16,17,18,19,21,22,24,25,26,27,32,33)

foreach ($r in $rowsToNo) {
    $ws.Cells.Item($r, 3).Value = "No"
}

# Update the selected range shown when the sheet is reopened.
$ws.Range("C31").Select()
